$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'69.794.34"
$ws.Range('E2').Value = '  +2.45%  '

$ws.Range('D3').Value = "'3.824.58"
$ws.Range('E3').Value = '  +0.84%  '

$ws.Range('D4').Value = "'0.999"
$ws.Range('E4').Value = '  -0.07%  '

$ws.Range('D5').Value = "'635.40"
$ws.Range('E5').Value = '  +5.76%  '

$ws.Range('D6').Value = "'165.84"
$ws.Range('E6').Value = '  +0.23%  '

$ws.Range('D7').Value = "'3.823.89"
$ws.Range('E7').Value = '  +0.94%  '

$ws.Range('E8').Value = '  -0.02%  '

$ws.Range('E9').Value = '  +0.91%  '

$ws.Range('D10').Value = "'0.162"
$ws.Range('E10').Value = '  +1.93%  '

$ws.Range('D11').Value = "'0.455"
$ws.Range('E11').Value = '  +0.78%  '

$ws.Range('D12').Value = "'6.72"
$ws.Range('E12').Value = '  +3.49%  '

$ws.Range('D13').Value = "'0.0000251"
$ws.Range('E13').Value = '  +0.52%  '

$ws.Range('D14').Value = "'36.01"
$ws.Range('E14').Value = '  +0.68%  '

$ws.Range('D15').Value = "'4.469.84"
$ws.Range('E15').Value = '  +0.90%  '

$ws.Range('D16').Value = "'3.821.12"
$ws.Range('E16').Value = '  +0.01%  '

$ws.Range('D17').Value = "'69.711.83"
$ws.Range('E17').Value = '  +2.34%  '

$ws.Range('D18').Value = "'18.09"
$ws.Range('E18').Value = '  -2.03%  '

$ws.Range('E19').Value = '  +1.18%  '

$ws.Range('E20').Value = '  -0.63%  '

$ws.Range('D21').Value = "'469.42"
$ws.Range('E21').Value = '  +1.88%  '

$ws.Range('D22').Value = "'9.73"
$ws.Range('E22').Value = '  +0.22%  '

$ws.Range('E23').Value = '  +1.50%  '

$ws.Range('D24').Value = "'0.0000151"
$ws.Range('E24').Value = '  +1.60%  '

$ws.Range('D25').Value = "'83.69"
$ws.Range('E25').Value = '  +0.99%  '

$ws.Range('B26').Value = 'InternetComputer(DFINITY)'
$ws.Range('C26').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D26').Value = "'12.12"
$ws.Range('E26').Value = '  +0.51%  '

$ws.Range('B27').Value = 'Fetch.AI'
$ws.Range('C27').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D27').Value = "'2.18"
$ws.Range('E27').Value = '  +3.25%  '

$ws.Range('D28').Value = "'10.08"
$ws.Range('E28').Value = '  +0.92%  '

$ws.Range('E29').Value = '  +0.00%  '

$ws.Range('D30').Value = "'3.978.37"
$ws.Range('E30').Value = '  +0.91%  '

$ws.Range('E31').Value = '  +2.46%  '

$ws.Range('B32').Value = 'ImmutableX'
$ws.Range('C32').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D32').Value = "'2.23"
$ws.Range('E32').Value = '  +0.38%  '

$ws.Range('B33').Value = 'NEARProtocol'
$ws.Range('C33').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D33').Value = "'7.33"
$ws.Range('E33').Value = '  -0.49%  '

$ws.Range('D34').Value = "'29.29"
$ws.Range('E34').Value = '  -0.12%  '

$ws.Range('B35').Value = 'Binance-PegBSC-USD'
$ws.Range('C35').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D35').Value = "'0.999"
$ws.Range('E35').Value = '  -0.04%  '

$ws.Range('B36').Value = 'RenzoRestakedETH'
$ws.Range('C36').Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range('D36').Value = "'3.765.66"
$ws.Range('E36').Value = '  +0.53%  '

$ws.Range('B37').Value = 'Aptos'
$ws.Range('C37').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D37').Value = "'9.08"
$ws.Range('E37').Value = '  +0.88%  '

$ws.Range('E38').Value = '  +3.32%  '

$ws.Range('E39').Value = '  +9.13%  '

$ws.Range('D40').Value = "'3.39"
$ws.Range('E40').Value = '  +3.10%  '

$ws.Range('D41').Value = "'5.91"
$ws.Range('E41').Value = '  +2.11%  '

$ws.Range('D42').Value = "'0.978"
$ws.Range('E42').Value = '  -0.72%  '

$ws.Range('E43').Value = '  -0.02%  '

$ws.Range('D45').Value = "'155.72"
$ws.Range('E45').Value = '  +2.66%  '

$ws.Range('D46').Value = "'0.301"
$ws.Range('E46').Value = '  +0.53%  '

$ws.Range('B47').Value = 'Stacks'
$ws.Range('C47').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D47').Value = "'1.95"
$ws.Range('E47').Value = '  +3.77%  '

$ws.Range('B48').Value = 'OKB'
$ws.Range('C48').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D48').Value = "'47.34"
$ws.Range('E48').Value = '  -0.05%  '

$ws.Range('B49').Value = 'Arweave'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D49').Value = "'43.28"
$ws.Range('E49').Value = '  -0.75%  '

$ws.Range('D50').Value = "'1.41"
$ws.Range('E50').Value = '  +3.81%  '

$ws.Range('E51').Value = '  +1.37%  '
